$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.744.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.942.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.32"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.507"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.940.43"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.431.20"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.730.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.942.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "441.50"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.41"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.666"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.32"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.16%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000102"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +16.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.44"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.989"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.14"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.59"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.64"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.03"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.48"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.280"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.52"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.697.84"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0337"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "362.47"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.83"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.97%  "
